$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final desired order for A2:C19 (player, position, team) - a reordering
# of the existing roster rows.
$data = @(
    @("Tyrese Maxey", "PG,SG", "Philadelphia 76ers"),
    @("Jaden Hardy", "PG,SG", "Dallas Mavericks"),
    @("Austin Reaves", "PG,SG", "Los Angeles Lakers"),
    @("Stephen Curry", "PG,SG", "Golden State Warriors"),
    @("OG Anunoby", "SF,PF", "New York Knicks"),
    @("Daniel Gafford", "PF,C", "Dallas Mavericks"),
    @("Karl-Anthony Towns", "PF,C", "New York Knicks"),
    @("Kevin Durant", "SF,PF", "Phoenix Suns"),
    @("Jarrett Allen", "C", "Cleveland Cavaliers"),
    @("Dereck Lively II", "C", "Dallas Mavericks"),
    @("Mark Williams", "C", "Charlotte Hornets"),
    @("Trey Murphy III", "SF,PF", "New Orleans Pelicans"),
    @("Darius Garland", "PG", "Cleveland Cavaliers"),
    @("Keegan Murray", "SF,PF", "Sacramento Kings"),
    @("Jalen Duren", "C", "Detroit Pistons"),
    @("Franz Wagner", "SF,PF", "Orlando Magic"),
    @("Jalen Johnson", "SF,PF", "Atlanta Hawks"),
    @("Tyrese Haliburton", "PG,SG", "Indiana Pacers")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
